# Add feature to refresh monthly cash balance.
# Rename the "kategorikas_listjumlahkas" sheet to "kategorikas_listsaldokas"
# (now showing a running/monthly balance instead of a sum) and make it the
# active/selected sheet, replacing "jenistransaksikas" as the active tab.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("kategorikas_listjumlahkas")
$ws.Name = "kategorikas_listsaldokas"
$ws.Activate()
